$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45696.42708333334, 96.61),
    @(45700.83333333334, 100),
    @(45701.73194444444, 75.16),
    @(45702.55833333333, 93.03),
    @(45702.59236111111, 87.34999999999999),
    @(45704.50138888889, 100),
    @(45707.08888888889, 88.95),
    @(45707.70208333333, 76.01000000000001),
    @(45710.84583333333, 97.15000000000001),
    @(45711.05902777778, 100)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
